$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "29.128.86"
Set-TextValue "E2" "  -0.39%  "
Set-TextValue "D3" "1.838.93"
Set-TextValue "E3" "  -0.32%  "
Set-TextValue "D4" "1.003"
Set-TextValue "E4" "  +0.43%  "
Set-TextValue "D5" "242.94"
Set-TextValue "E5" "  +0.12%  "
Set-TextValue "D6" "0.6250"
Set-TextValue "E6" "  -5.73%  "
Set-TextValue "D7" "1.005"
Set-TextValue "E7" "  +0.51%  "
Set-TextValue "D8" "0.07592"
Set-TextValue "E8" "  +1.99%  "
Set-TextValue "D9" "0.2922"
Set-TextValue "E9" "  -1.16%  "
Set-TextValue "D10" "22.69"
Set-TextValue "E10" "  -2.71%  "
Set-TextValue "D11" "0.07767"
Set-TextValue "E11" "  +0.06%  "
Set-TextValue "D12" "1.843.99"
Set-TextValue "E12" "  +0.07%  "
Set-TextValue "D13" "4.959"
Set-TextValue "D14" "0.6662"
Set-TextValue "E14" "  -0.96%  "
Set-TextValue "D15" "0.000009978"
Set-TextValue "E15" "  +14.48%  "
Set-TextValue "D16" "82.81"
Set-TextValue "E16" "  -0.82%  "
Set-TextValue "D17" "6.030"
Set-TextValue "E17" "  -2.54%  "
Set-TextValue "D18" "29.163.53"
Set-TextValue "E18" "  -0.27%  "
Set-TextValue "D19" "225.62"
Set-TextValue "E19" "  -0.53%  "
Set-TextValue "D20" "12.36"
Set-TextValue "E20" "  -1.45%  "
Set-TextValue "D21" "1.005"
Set-TextValue "E21" "  +0.44%  "
Set-TextValue "D22" "7.203"
Set-TextValue "D23" "1.005"
Set-TextValue "E23" "  +0.54%  "
Set-TextValue "D24" "159.04"
Set-TextValue "E24" "  +0.18%  "
Set-TextValue "D25" "8.470"
Set-TextValue "E25" "  -1.96%  "
Set-TextValue "D26" "0.1366"
Set-TextValue "E26" "  -2.94%  "
Set-TextValue "D27" "17.94"
Set-TextValue "E27" "  -0.56%  "
Set-TextValue "D28" "1.491"
Set-TextValue "E28" "  -1.14%  "
Set-TextValue "D29" "4.077"
Set-TextValue "D30" "4.040"
Set-TextValue "E30" "  -0.45%  "
Set-TextValue "D31" "1.202"
Set-TextValue "E31" "  +0.86%  "
Set-TextValue "D32" "0.05206"
Set-TextValue "E32" "  -2.35%  "
Set-TextValue "D33" "1.853"
Set-TextValue "E33" "  -1.03%  "
Set-TextValue "D34" "0.7388"
Set-TextValue "E34" "  -1.28%  "
Set-TextValue "D35" "1.143"
Set-TextValue "E35" "  -1.24%  "
Set-TextValue "D36" "2.712"
Set-TextValue "E36" "  +2.08%  "
Set-TextValue "D37" "1.254.43"
Set-TextValue "E37" "  -4.61%  "
Set-TextValue "D38" "2.771"
Set-TextValue "E38" "  +0.47%  "
Set-TextValue "D39" "0.01784"
Set-TextValue "E39" "  -1.08%  "
Set-TextValue "D40" "6.325"
Set-TextValue "E40" "  -1.26%  "
Set-TextValue "D41" "0.8968"
Set-TextValue "E41" "  -0.87%  "
Set-TextValue "D42" "1.006"
Set-TextValue "E42" "  +0.64%  "
Set-TextValue "D43" "101.49"
Set-TextValue "E43" "  -1.97%  "
Set-TextValue "B44" "RocketPoolETH"
Set-TextValue "C44" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D44" "1.979.46"
Set-TextValue "E44" "  -0.61%  "
Set-TextValue "B45" "BabyDogeCoin"
Set-TextValue "C45" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D45" "0.00000000123"
Set-TextValue "E45" "  +0.81%  "
Set-TextValue "D46" "64.16"
Set-TextValue "E46" "  -1.80%  "
Set-TextValue "D47" "0.5130"
Set-TextValue "E47" "  -0.26%  "
Set-TextValue "D48" "0.4016"
Set-TextValue "E48" "  -0.21%  "
Set-TextValue "D49" "8.857"
Set-TextValue "E49" "  +1.25%  "
Set-TextValue "D50" "0.05760"
Set-TextValue "E50" "  -1.82%  "
Set-TextValue "D51" "1.644"
Set-TextValue "E51" "  -6.29%  "
